# Refresh NATMI Reln -> Itga3 LR-pair TPM stats (adds "Resolving-Mac" as a 4th
# sending/target cluster, rows 2-17) to reflect the updated TPM recomputation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Reln"
$ws.Cells.Item(2, 3).Value = "Itga3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.06455366666666666
$ws.Cells.Item(2, 8).Value = 0.193661
$ws.Cells.Item(2, 9).Value = 0.01357839286814829
$ws.Cells.Item(2, 10).Value = 0.01357839286814829
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 13.604331
$ws.Cells.Item(2, 14).Value = 40.812993
$ws.Cells.Item(2, 15).Value = 0.8107276168878804
$ws.Cells.Item(2, 16).Value = 0.8107276168878805
$ws.Cells.Item(2, 17).Value = 0.8782094485969999
$ws.Cells.Item(2, 18).Value = 7.903885037373
$ws.Cells.Item(2, 19).Value = 0.01100837809116125
$ws.Cells.Item(2, 20).Value = 0.01100837809116125

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Reln"
$ws.Cells.Item(3, 3).Value = "Itga3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.06455366666666666
$ws.Cells.Item(3, 8).Value = 0.193661
$ws.Cells.Item(3, 9).Value = 0.01357839286814829
$ws.Cells.Item(3, 10).Value = 0.01357839286814829
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.061748
$ws.Cells.Item(3, 14).Value = 3.185244
$ws.Cells.Item(3, 15).Value = 0.06327311690486458
$ws.Cells.Item(3, 16).Value = 0.06327311690486459
$ws.Cells.Item(3, 17).Value = 0.06853972647599998
$ws.Cells.Item(3, 18).Value = 0.616857538284
$ws.Cells.Item(3, 19).Value = 0.000859147239326526
$ws.Cells.Item(3, 20).Value = 0.0008591472393265264

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Reln"
$ws.Cells.Item(4, 3).Value = "Itga3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.06455366666666666
$ws.Cells.Item(4, 8).Value = 0.193661
$ws.Cells.Item(4, 9).Value = 0.01357839286814829
$ws.Cells.Item(4, 10).Value = 0.01357839286814829
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.995771333333333
$ws.Cells.Item(4, 14).Value = 5.987314
$ws.Cells.Item(4, 15).Value = 0.1189346934389115
$ws.Cells.Item(4, 16).Value = 0.1189346934389116
$ws.Cells.Item(4, 17).Value = 0.1288343573948889
$ws.Cells.Item(4, 18).Value = 1.159509216554
$ws.Cells.Item(4, 19).Value = 0.001614941993166319
$ws.Cells.Item(4, 20).Value = 0.00161494199316632

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Reln"
$ws.Cells.Item(5, 3).Value = "Itga3"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.06455366666666666
$ws.Cells.Item(5, 8).Value = 0.193661
$ws.Cells.Item(5, 9).Value = 0.01357839286814829
$ws.Cells.Item(5, 10).Value = 0.01357839286814829
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1185463333333334
$ws.Cells.Item(5, 14).Value = 0.355639
$ws.Cells.Item(5, 15).Value = 0.007064572768343379
$ws.Cells.Item(5, 16).Value = 0.007064572768343379
$ws.Cells.Item(5, 17).Value = 0.007652600486555556
$ws.Cells.Item(5, 18).Value = 0.06887340437900001
$ws.Cells.Item(5, 19).Value = 0.00009592554449418834
$ws.Cells.Item(5, 20).Value = 0.00009592554449418836

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Reln"
$ws.Cells.Item(6, 3).Value = "Itga3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.06084700000000001
$ws.Cells.Item(6, 8).Value = 0.182541
$ws.Cells.Item(6, 9).Value = 0.01279872257472933
$ws.Cells.Item(6, 10).Value = 0.01279872257472933
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 13.604331
$ws.Cells.Item(6, 14).Value = 40.812993
$ws.Cells.Item(6, 15).Value = 0.8107276168878804
$ws.Cells.Item(6, 16).Value = 0.8107276168878805
$ws.Cells.Item(6, 17).Value = 0.8277827283570001
$ws.Cells.Item(6, 18).Value = 7.450044555213
$ws.Cells.Item(6, 19).Value = 0.01037627785221943
$ws.Cells.Item(6, 20).Value = 0.01037627785221943

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Reln"
$ws.Cells.Item(7, 3).Value = "Itga3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.06084700000000001
$ws.Cells.Item(7, 8).Value = 0.182541
$ws.Cells.Item(7, 9).Value = 0.01279872257472933
$ws.Cells.Item(7, 10).Value = 0.01279872257472933
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.061748
$ws.Cells.Item(7, 14).Value = 3.185244
$ws.Cells.Item(7, 15).Value = 0.06327311690486458
$ws.Cells.Item(7, 16).Value = 0.06327311690486459
$ws.Cells.Item(7, 17).Value = 0.06460418055600001
$ws.Cells.Item(7, 18).Value = 0.581437625004
$ws.Cells.Item(7, 19).Value = 0.0008098150697037784
$ws.Cells.Item(7, 20).Value = 0.0008098150697037784

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Reln"
$ws.Cells.Item(8, 3).Value = "Itga3"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.06084700000000001
$ws.Cells.Item(8, 8).Value = 0.182541
$ws.Cells.Item(8, 9).Value = 0.01279872257472933
$ws.Cells.Item(8, 10).Value = 0.01279872257472933
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.995771333333333
$ws.Cells.Item(8, 14).Value = 5.987314
$ws.Cells.Item(8, 15).Value = 0.1189346934389115
$ws.Cells.Item(8, 16).Value = 0.1189346934389116
$ws.Cells.Item(8, 17).Value = 0.1214366983193333
$ws.Cells.Item(8, 18).Value = 1.092930284874
$ws.Cells.Item(8, 19).Value = 0.00152221214583511
$ws.Cells.Item(8, 20).Value = 0.00152221214583511

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Reln"
$ws.Cells.Item(9, 3).Value = "Itga3"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.06084700000000001
$ws.Cells.Item(9, 8).Value = 0.182541
$ws.Cells.Item(9, 9).Value = 0.01279872257472933
$ws.Cells.Item(9, 10).Value = 0.01279872257472933
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1185463333333334
$ws.Cells.Item(9, 14).Value = 0.355639
$ws.Cells.Item(9, 15).Value = 0.007064572768343379
$ws.Cells.Item(9, 16).Value = 0.007064572768343379
$ws.Cells.Item(9, 17).Value = 0.007213188744333335
$ws.Cells.Item(9, 18).Value = 0.06491869869900001
$ws.Cells.Item(9, 19).Value = 0.0000904175069710145
$ws.Cells.Item(9, 20).Value = 0.00009041750697101448

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Reln"
$ws.Cells.Item(10, 3).Value = "Itga3"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.62452
$ws.Cells.Item(10, 8).Value = 13.87356
$ws.Cells.Item(10, 9).Value = 0.9727340463997778
$ws.Cells.Item(10, 10).Value = 0.9727340463997778
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.604331
$ws.Cells.Item(10, 14).Value = 40.812993
$ws.Cells.Item(10, 15).Value = 0.8107276168878804
$ws.Cells.Item(10, 16).Value = 0.8107276168878805
$ws.Cells.Item(10, 17).Value = 62.91350079612
$ws.Cells.Item(10, 18).Value = 566.22150716508
$ws.Cells.Item(10, 19).Value = 0.7886223553033967
$ws.Cells.Item(10, 20).Value = 0.7886223553033969

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Reln"
$ws.Cells.Item(11, 3).Value = "Itga3"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.62452
$ws.Cells.Item(11, 8).Value = 13.87356
$ws.Cells.Item(11, 9).Value = 0.9727340463997778
$ws.Cells.Item(11, 10).Value = 0.9727340463997778
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.061748
$ws.Cells.Item(11, 14).Value = 3.185244
$ws.Cells.Item(11, 15).Value = 0.06327311690486458
$ws.Cells.Item(11, 16).Value = 0.06327311690486459
$ws.Cells.Item(11, 17).Value = 4.910074860959999
$ws.Cells.Item(11, 18).Value = 44.19067374863999
$ws.Cells.Item(11, 19).Value = 0.06154791503519511
$ws.Cells.Item(11, 20).Value = 0.06154791503519512

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Reln"
$ws.Cells.Item(12, 3).Value = "Itga3"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.62452
$ws.Cells.Item(12, 8).Value = 13.87356
$ws.Cells.Item(12, 9).Value = 0.9727340463997778
$ws.Cells.Item(12, 10).Value = 0.9727340463997778
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.995771333333333
$ws.Cells.Item(12, 14).Value = 5.987314
$ws.Cells.Item(12, 15).Value = 0.1189346934389115
$ws.Cells.Item(12, 16).Value = 0.1189346934389116
$ws.Cells.Item(12, 17).Value = 9.229484446426664
$ws.Cells.Item(12, 18).Value = 83.06536001783999
$ws.Cells.Item(12, 19).Value = 0.1156918256061495
$ws.Cells.Item(12, 20).Value = 0.1156918256061495

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Reln"
$ws.Cells.Item(13, 3).Value = "Itga3"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.62452
$ws.Cells.Item(13, 8).Value = 13.87356
$ws.Cells.Item(13, 9).Value = 0.9727340463997778
$ws.Cells.Item(13, 10).Value = 0.9727340463997778
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1185463333333334
$ws.Cells.Item(13, 14).Value = 0.355639
$ws.Cells.Item(13, 15).Value = 0.007064572768343379
$ws.Cells.Item(13, 16).Value = 0.007064572768343379
$ws.Cells.Item(13, 17).Value = 0.5482198894266667
$ws.Cells.Item(13, 18).Value = 4.93397900484
$ws.Cells.Item(13, 19).Value = 0.006871950455036335
$ws.Cells.Item(13, 20).Value = 0.006871950455036335

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Reln"
$ws.Cells.Item(14, 3).Value = "Itga3"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.004225666666666667
$ws.Cells.Item(14, 8).Value = 0.012677
$ws.Cells.Item(14, 9).Value = 0.0008888381573446169
$ws.Cells.Item(14, 10).Value = 0.000888838157344617
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 13.604331
$ws.Cells.Item(14, 14).Value = 40.812993
$ws.Cells.Item(14, 15).Value = 0.8107276168878804
$ws.Cells.Item(14, 16).Value = 0.8107276168878805
$ws.Cells.Item(14, 17).Value = 0.057487368029
$ws.Cells.Item(14, 18).Value = 0.5173863122610001
$ws.Cells.Item(14, 19).Value = 0.0007206056411030162
$ws.Cells.Item(14, 20).Value = 0.0007206056411030163

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Reln"
$ws.Cells.Item(15, 3).Value = "Itga3"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.004225666666666667
$ws.Cells.Item(15, 8).Value = 0.012677
$ws.Cells.Item(15, 9).Value = 0.0008888381573446169
$ws.Cells.Item(15, 10).Value = 0.000888838157344617
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.061748
$ws.Cells.Item(15, 14).Value = 3.185244
$ws.Cells.Item(15, 15).Value = 0.06327311690486458
$ws.Cells.Item(15, 16).Value = 0.06327311690486459
$ws.Cells.Item(15, 17).Value = 0.004486593132
$ws.Cells.Item(15, 18).Value = 0.04037933818800001
$ws.Cells.Item(15, 19).Value = 0.00005623956063917036
$ws.Cells.Item(15, 20).Value = 0.00005623956063917038

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Reln"
$ws.Cells.Item(16, 3).Value = "Itga3"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.004225666666666667
$ws.Cells.Item(16, 8).Value = 0.012677
$ws.Cells.Item(16, 9).Value = 0.0008888381573446169
$ws.Cells.Item(16, 10).Value = 0.000888838157344617
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.995771333333333
$ws.Cells.Item(16, 14).Value = 5.987314
$ws.Cells.Item(16, 15).Value = 0.1189346934389115
$ws.Cells.Item(16, 16).Value = 0.1189346934389116
$ws.Cells.Item(16, 17).Value = 0.008433464397555554
$ws.Cells.Item(16, 18).Value = 0.075901179578
$ws.Cells.Item(16, 19).Value = 0.000105713693760589
$ws.Cells.Item(16, 20).Value = 0.0001057136937605891

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Reln"
$ws.Cells.Item(17, 3).Value = "Itga3"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.004225666666666667
$ws.Cells.Item(17, 8).Value = 0.012677
$ws.Cells.Item(17, 9).Value = 0.0008888381573446169
$ws.Cells.Item(17, 10).Value = 0.000888838157344617
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1185463333333334
$ws.Cells.Item(17, 14).Value = 0.355639
$ws.Cells.Item(17, 15).Value = 0.007064572768343379
$ws.Cells.Item(17, 16).Value = 0.007064572768343379
$ws.Cells.Item(17, 17).Value = 0.0005009372892222223
$ws.Cells.Item(17, 18).Value = 0.004508435603
$ws.Cells.Item(17, 19).Value = 0.000006279261841841289
$ws.Cells.Item(17, 20).Value = 0.000006279261841841289

